$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last header cell (G1) onto the new H1 header
# cell so the new "Save" column matches the existing header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data value for the new "Save" column (row 2), left unstyled like
# the other numeric cells in that row.
$ws.Range("H2").Value = 1
